# Add experimental data on 2021-10-29
# Also: rename "Shapes_N" subfolders to "Group_N" for previously-entered
# experiments, and switch the "Exp date" column to free-text dd/mm/yyyy
# strings (row 7 keeps a raw date serial, matching the source edit).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename existing "Shapes_N" paths to "Group_N" (col B) and switch
#     "Exp date" (col A) for rows 2-6 to text dd/mm/yyyy -------------------
$ws.Range("A2").Value = "13/04/2021"
$ws.Range("B2").Value = "F:\PhD, PMMH, ESPCI\Processing\20210413-Actin\results\Group_1"

$ws.Range("A3").Value = "30/04/2021"
$ws.Range("B3").Value = "F:\PhD, PMMH, ESPCI\Processing\20210430-Actin\results\Group_1"

$ws.Range("A4").Value = "30/04/2021"
$ws.Range("B4").Value = "F:\PhD, PMMH, ESPCI\Processing\20210430-Actin\results\Group_2"

$ws.Range("A5").Value = "25/09/2020"
$ws.Range("B5").Value = "F:\PhD, PMMH, ESPCI\Processing\20200925-Actin\results\Group_1"

$ws.Range("A6").Value = "25/09/2020"
$ws.Range("B6").Value = "F:\PhD, PMMH, ESPCI\Processing\20200925-Actin\results\Group_2"

# Row 7 keeps a numeric date, but the value itself changes.
$ws.Range("A7").Value = 43840
$ws.Range("B7").Value = "F:\PhD, PMMH, ESPCI\Processing\20201001-Actin\results\Group_1"

# --- New experiment rows for 2021-10-29 -----------------------------------
$ws.Range("A8").Value = "29/10/2021"
$ws.Range("B8").Value = "G:\PhD, PMMH, ESPCI\Processing\20211029-Actin\results\Group_1"
$ws.Range("C8").Value = "G:\PhD, PMMH, ESPCI\Processing\20211029-Actin\results\circlesforPAs20211029_S15.mat"
$ws.Range("D8").Value = "G:\PhD, PMMH, ESPCI\Processing\20211029-Actin\results\Figures"
$ws.Range("E8").Value = 55
$ws.Range("F8").Value = 400
$ws.Range("G8").Value = 0.5
$ws.Range("H8").Formula = "=G8/F8/E8"
$ws.Range("I8").Value = 0.1
$ws.Range("J8").Value = 20

$ws.Range("A9").Value = "29/10/2021"
$ws.Range("B9").Value = "G:\PhD, PMMH, ESPCI\Processing\20211029-Actin\results\Group_2"
$ws.Range("C9").Value = "G:\PhD, PMMH, ESPCI\Processing\20211029-Actin\results\circlesforPAs20211029_S15.mat"
$ws.Range("D9").Value = "G:\PhD, PMMH, ESPCI\Processing\20211029-Actin\results\Figures"
$ws.Range("E9").Value = 55
$ws.Range("F9").Value = 400
$ws.Range("G9").Value = 1
$ws.Range("H9").Formula = "=G9/F9/E9"
$ws.Range("I9").Value = 0.1
$ws.Range("J9").Value = 20

$ws.Range("A10").Value = "29/10/2021"
$ws.Range("B10").Value = "G:\PhD, PMMH, ESPCI\Processing\20211029-Actin\results\Group_3"
$ws.Range("C10").Value = "G:\PhD, PMMH, ESPCI\Processing\20211029-Actin\results\circlesforPAs20211029_S15.mat"
$ws.Range("D10").Value = "G:\PhD, PMMH, ESPCI\Processing\20211029-Actin\results\Figures"
$ws.Range("E10").Value = 55
$ws.Range("F10").Value = 400
$ws.Range("G10").Value = 1.5
$ws.Range("H10").Formula = "=G10/F10/E10"
$ws.Range("I10").Value = 0.1
$ws.Range("J10").Value = 20

# --- Column C needs to grow to fit the new (longer) path strings ---------
$ws.Columns.Item(3).ColumnWidth = 81.1

# --- Update selection to match the source edit -----------------------------
$ws.Range("E13").Select()
